$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: set NumberFormat to text ("@") before writing price strings to
# the Price column so Excel does not silently coerce them into numeric values
# (these source values are plain text, e.g. "1.867.37" or "0.000007455"),
# then restore the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "30.019.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  +9.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.868.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  +6.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "250.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4973"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  +3.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "45.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +9.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.2837"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  +8.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.06527"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +5.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "1.870.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +7.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "16.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +5.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.07209"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +3.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.6592"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +9.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "85.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +9.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "4.803"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +7.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "29.998.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +9.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "2.109.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +7.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "4.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +6.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "9.014"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +6.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "5.487"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +7.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "144.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "134.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +24.32%  "
$ws.Range("E28").Value2 = "  +9.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.939"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +5.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.395"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "4.233"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +7.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.08572"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +7.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.875"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +5.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.05058"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +7.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +11.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.6812"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +10.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "2.693"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +3.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.330"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +15.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "2.731"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +6.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.9585"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +3.60%  "
$ws.Range("E41").Value2 = "  +8.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "6.131"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +7.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "103.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.4162"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +8.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "7.441"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +7.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.1249"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +7.91%  "
$ws.Range("E48").Value2 = "  +5.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "8.286"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +5.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "32.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +8.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.3705"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +9.68%  "

# Row 19 (was ShibaInu) becomes Avalanche
$ws.Range("B19").Value2 = "Avalanche"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "12.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +11.26%  "

# Row 20 (was Avalanche) becomes ShibaInu
$ws.Range("B20").Value2 = "ShibaInu"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "0.000007467"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +5.71%  "
